$wb = $excel.ActiveWorkbook
$ov = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This shared string is used by Overview!E2, Overview!F2, zh-cn!C2 and de-de!C2.
$ov.Range("E2").Value = "Handed back: in sync with en-US"
$ov.Range("F2").Value = "Handed back: in sync with en-US"
$zh.Range("C2").Value = "Handed back: in sync with en-US"
$de.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn row 2: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$zh.Range("I2").Value = "47181a61-8ac6-4af1-b013-0f00d62e9ca2.md"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e90a835daa426dec55fcd05ab54005cb0c2f3b8e/e2e/47181a61-8ac6-4af1-b013-0f00d62e9ca2.md", $null, $null, "47181a61-8ac6-4af1-b013-0f00d62e9ca2.md")
$zh.Range("I2").Font.Name = "Calibri"
$zh.Range("I2").Font.Size = 11
$zh.Range("I2").Font.Underline = 2
$zh.Range("I2").Font.Color = 15570276

$zh.Range("J2").Value = "47181a61-8ac6-4af1-b013-0f00d62e9ca2.a4442fc7ddc6ac07660181dde5277a26e67b1425.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-12 21:12:29"

# --- de-de row 2: Latest Target File / Latest Handback File / Latest Handback DateTime ---
$de.Range("I2").Value = "47181a61-8ac6-4af1-b013-0f00d62e9ca2.md"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/oltest/blob/e90a835daa426dec55fcd05ab54005cb0c2f3b8e/e2e/47181a61-8ac6-4af1-b013-0f00d62e9ca2.md", $null, $null, "47181a61-8ac6-4af1-b013-0f00d62e9ca2.md")
$de.Range("I2").Font.Name = "Calibri"
$de.Range("I2").Font.Size = 11
$de.Range("I2").Font.Underline = 2
$de.Range("I2").Font.Color = 15570276

$de.Range("J2").Value = "47181a61-8ac6-4af1-b013-0f00d62e9ca2.a4442fc7ddc6ac07660181dde5277a26e67b1425.de-de.xlf"
$de.Range("K2").Value = "2016-08-12 21:12:38"

# --- Column width adjustments (autofit-style widening to fit new longer text) ---
$ov.Range("E1").EntireColumn.ColumnWidth = 29.166666666666668
$ov.Range("F1").EntireColumn.ColumnWidth = 29.166666666666668

$zh.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$zh.Range("I1").EntireColumn.ColumnWidth = 39.166666666666664
$zh.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664

$de.Range("C1").EntireColumn.ColumnWidth = 29.166666666666668
$de.Range("I1").EntireColumn.ColumnWidth = 39.166666666666664
$de.Range("J1").EntireColumn.ColumnWidth = 39.166666666666664
